# Append new conversation rows (16-18) to the active worksheet,
# matching the rows already present (Hora, Sujeto, Mensaje columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "2025-09-23 11:29:21"
$ws.Range("B16").Value = "Usuario"
$ws.Range("C16").Value = "hola"

$ws.Range("A17").Value = "2025-09-23 11:33:11"
$ws.Range("B17").Value = "Usuario"
$ws.Range("C17").Value = "hola"

$ws.Range("A18").Value = "2025-09-23 11:33:15"
$ws.Range("B18").Value = "Asistente"
$ws.Range("C18").Value = "¡Hola! Soy Seraphina, tu asistente de bienestar integral. ¿En qué área de tu bienestar te gustaría enfocarte hoy?
"
